# Adapt to new ver
# - SONG_ID (A2): dogma1 -> anochu
# - GENRE (B2): 0
# - TITLE (C2): "ちゅ、多様性。"
# - SUBTITLE (D2): "TVアニメ「チェンソーマン」エンディング・テーマ"
# - ONI_LEVEL (F2): 8
# - HARD_LEVEL (G2): 5
# - NORMAL_LEVEL (H2): 4
# - EASY_LEVEL (I2): 3
# - move selection from A2 to A3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "anochu"
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = "ちゅ、多様性。"
$ws.Range("D2").Value = "TVアニメ「チェンソーマン」エンディング・テーマ"
$ws.Range("F2").Value = 8
$ws.Range("G2").Value = 5
$ws.Range("H2").Value = 4
$ws.Range("I2").Value = 3

$ws.Range("A3").Select()
